$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 previously held only an empty, bordered "footer" pair of cells (B13/C13).
# Turn it into a real data row - "Step4" / its formula text - matching the
# pattern already used by rows 10-12 (plain label in column B, quote-prefixed
# formula-looking text in column C), and drop the border that used to mark the
# bottom of the table since it is no longer the last row.
$ws.Range("B13:C13").Borders.LineStyle = 0

$ws.Range("B13").Value = "Step4"
$ws.Range("C13").Value = "'= OpenLUtils.dateToString(new Date(0))"

[void]$ws.Range("C13").Select()
